# Auto-generated data refresh: updates computed market-price/profit columns
# (H..N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve tables to match
# the latest scheduled-runner snapshot. All target cells are plain numeric
# literals (no formulas anywhere in this workbook).
$wb = $excel.ActiveWorkbook

# ALC row 6: "Days of Chunder" / "Antidote"
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 44969.6
$ws.Range("J6").Value = 56158.125
$ws.Range("L6").Value = 168474.375
$ws.Range("N6").Value = -168698.375

# ALC row 86: "Filling in the Blanks" / "Enchanted Aurum Regis Ink"
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2499.3333
$ws.Range("I86").Value = 1990
$ws.Range("J86").Value = 2601.2
$ws.Range("K86").Value = 1990
$ws.Range("L86").Value = 2601.2
$ws.Range("M86").Value = -867
$ws.Range("N86").Value = -4847.2

# ALC row 88: "The Grave of Hemlock Groves" / "Growth Formula Zeta"
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 22209.5
$ws.Range("I88").Value = 2301.5
$ws.Range("J88").Value = 27186.5
$ws.Range("K88").Value = 2301.5
$ws.Range("L88").Value = 27186.5
$ws.Range("M88").Value = -1895.5
$ws.Range("N88").Value = -27998.5

# ALC row 89: "Ink into Antiquity (L)" / "Enchanted Aurum Regis Ink"
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 2499.3333
$ws.Range("I89").Value = 1990
$ws.Range("J89").Value = 2601.2
$ws.Range("K89").Value = 9950
$ws.Range("L89").Value = 13006
$ws.Range("M89").Value = -4334
$ws.Range("N89").Value = -24238

# ALC row 91: "Dappling the Highlands (L)" / "Growth Formula Zeta"
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 22209.5
$ws.Range("I91").Value = 2301.5
$ws.Range("J91").Value = 27186.5
$ws.Range("K91").Value = 2301.5
$ws.Range("L91").Value = 27186.5
$ws.Range("M91").Value = -897.5
$ws.Range("N91").Value = -29994.5

# ALC row 98: "The Dotted Line" / "Enchanted Durium Ink"
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2076.7827
$ws.Range("I98").Value = 1157.5
$ws.Range("J98").Value = 2567.0667
$ws.Range("K98").Value = 1157.5
$ws.Range("L98").Value = 2567.0667
$ws.Range("M98").Value = 340.5
$ws.Range("N98").Value = -5563.066699999999

# ALC row 107: "Another Man's Ink" / "Enchanted Truegold Ink"
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1653
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 1653
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -5493

# ALC row 116: "Growing Up" / "Growth Formula Kappa"
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 1777.6666
$ws.Range("I116").Value = 1741.5
$ws.Range("J116").Value = 1850
$ws.Range("K116").Value = 1741.5
$ws.Range("L116").Value = 1850
$ws.Range("M116").Value = 1700.5
$ws.Range("N116").Value = -8734

# ALC row 122: "Wishful Inking" / "Enchanted High Durium Ink"
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2076.7827
$ws.Range("I122").Value = 1157.5
$ws.Range("J122").Value = 2567.0667
$ws.Range("K122").Value = 3472.5
$ws.Range("L122").Value = 7701.2001
$ws.Range("M122").Value = -1022.5
$ws.Range("N122").Value = -12601.2001

# ARM row 5: "The Alloyed Truth" / "Bronze Rivets"
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 17543952
$ws.Range("I5").Value = 17543952
$ws.Range("K5").Value = 17543952
$ws.Range("M5").Value = -17543840

# BSM row 4: "Mending Fences" / "Bronze Rivets"
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 17543952
$ws.Range("I4").Value = 17543952
$ws.Range("K4").Value = 17543952
$ws.Range("M4").Value = -17543837

# BSM row 80: "Unbreaker" / "Titanium Ingot"
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 502.6
$ws.Range("I80").Value = 100.666664
$ws.Range("J80").Value = 573.5294
$ws.Range("K80").Value = 100.666664
$ws.Range("L80").Value = 573.5294
$ws.Range("M80").Value = 897.333336
$ws.Range("N80").Value = -2569.5294

# BSM row 83: "Attack on Titanium (L)" / "Titanium Ingot"
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 502.6
$ws.Range("I83").Value = 100.666664
$ws.Range("J83").Value = 573.5294
$ws.Range("K83").Value = 503.33332
$ws.Range("L83").Value = 2867.647
$ws.Range("M83").Value = 4488.66668
$ws.Range("N83").Value = -12851.647

# CRP row 58: "You Do the Heavy Lifting" / "Mahogany Lumber"
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 35716360
$ws.Range("I58").Value = 83334510
$ws.Range("J58").Value = 2748.75
$ws.Range("K58").Value = 83334510
$ws.Range("L58").Value = 2748.75
$ws.Range("M58").Value = -83334307
$ws.Range("N58").Value = -3154.75

# CRP row 99: "O Pine" / "Pine Lumber"
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3339.1365
$ws.Range("I99").Value = 3428.3572
$ws.Range("J99").Value = 3183
$ws.Range("K99").Value = 3428.3572
$ws.Range("L99").Value = 3183
$ws.Range("M99").Value = -1930.3572
$ws.Range("N99").Value = -6179

# CRP row 107: "Built to Last" / "White Oak Lumber"
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 368.1111
$ws.Range("I107").Value = 219.5625
$ws.Range("J107").Value = 1556.5
$ws.Range("K107").Value = 219.5625
$ws.Range("L107").Value = 1556.5
$ws.Range("M107").Value = 1700.4375
$ws.Range("N107").Value = -5396.5

# CRP row 126: "A Better Conductor" / "Red Pine Lumber"
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3339.1365
$ws.Range("I126").Value = 3428.3572
$ws.Range("J126").Value = 3183
$ws.Range("K126").Value = 10285.0716
$ws.Range("L126").Value = 9549
$ws.Range("M126").Value = -7815.071599999999
$ws.Range("N126").Value = -14489

# CRP row 129: "Spinning the Time Away" / "Ironwood Spinning Wheel"
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H129").Value = 31723.5
$ws.Range("J129").Value = 31723.5
$ws.Range("L129").Value = 31723.5
$ws.Range("N129").Value = -41723.5

# CRP row 132: "Hull Lotta Damage" / "Ginseng Lumber"
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 94949.82000000001
$ws.Range("I132").Value = 4350
$ws.Range("J132").Value = 253499.5
$ws.Range("K132").Value = 13050
$ws.Range("L132").Value = 760498.5
$ws.Range("M132").Value = -10520
$ws.Range("N132").Value = -765558.5

# CRP row 134: "Wood You Be Quiet" / "Ceiba Lumber"
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 35561.906
$ws.Range("I134").Value = 2248.1
$ws.Range("K134").Value = 6744.299999999999
$ws.Range("M134").Value = -4209.299999999999

# CRP row 136: "Turali Quality" / "Dark Mahogany Lumber"
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 35716360
$ws.Range("I136").Value = 83334510
$ws.Range("J136").Value = 2748.75
$ws.Range("K136").Value = 250003530
$ws.Range("L136").Value = 8246.25
$ws.Range("M136").Value = -250000980
$ws.Range("N136").Value = -13346.25

# CUL row 17: "Chew the Fat" / "Grilled Dodo"
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 840
$ws.Range("I17").Value = 800
$ws.Range("K17").Value = 2400
$ws.Range("M17").Value = -2231

# CUL row 131: "The Mountain Steeped" / "Tsai tou Vounou"
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1312.6038
$ws.Range("I131").Value = 684.1667
$ws.Range("J131").Value = 1496.5366
$ws.Range("K131").Value = 2052.5001
$ws.Range("L131").Value = 4489.6098
$ws.Range("M131").Value = 2987.4999
$ws.Range("N131").Value = -14569.6098

# GSM row 7: "Water of Life" / "Copper Rings"
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 40000000
$ws.Range("I7").Value = 40000000
$ws.Range("K7").Value = 40000000
$ws.Range("M7").Value = -39999888

# GSM row 8: "Gods of Small Things" / "Copper Earrings"
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H8").Value = 40000000
$ws.Range("I8").Value = 40000000
$ws.Range("K8").Value = 40000000
$ws.Range("M8").Value = -39999861

# GSM row 122: "Awarding Academic Excellence" / "Ametrine"
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4745
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 5993.3335
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 17980.0005
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -22880.0005

# LTW row 7: "Tan Before the Ban" / "Leather"
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1999.8823
$ws.Range("I7").Value = 1807.8334
$ws.Range("J7").Value = 2460.8
$ws.Range("K7").Value = 1807.8334
$ws.Range("L7").Value = 2460.8
$ws.Range("M7").Value = -1695.8334
$ws.Range("N7").Value = -2684.8

# LTW row 40: "Best Served Toad" / "Toad Leather"
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3219.9
$ws.Range("I40").Value = 2813.2
$ws.Range("J40").Value = 4440
$ws.Range("K40").Value = 2813.2
$ws.Range("L40").Value = 4440
$ws.Range("M40").Value = -2677.2
$ws.Range("N40").Value = -4712

# LTW row 122: "Hell on Leather" / "Gaja Leather"
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2551.3
$ws.Range("I122").Value = 2495.375
$ws.Range("J122").Value = 2775
$ws.Range("K122").Value = 7486.125
$ws.Range("L122").Value = 8325
$ws.Range("M122").Value = -5036.125
$ws.Range("N122").Value = -13225

# LTW row 126: "Battered Books" / "Saiga Leather"
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1999.8823
$ws.Range("I126").Value = 1807.8334
$ws.Range("J126").Value = 2460.8
$ws.Range("K126").Value = 5423.5002
$ws.Range("L126").Value = 7382.400000000001
$ws.Range("M126").Value = -2953.5002
$ws.Range("N126").Value = -12322.4

# LTW row 132: "Tenets of Tanning" / "Silver Lobo Leather"
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 41969.8
$ws.Range("I132").Value = 1236.7778
$ws.Range("K132").Value = 3710.3334
$ws.Range("M132").Value = -1180.3334

# WVR row 96: "Skills on Display" / "Ruby Cotton Cloth"
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 6616.6665
$ws.Range("I96").Value = 5185.7144
$ws.Range("J96").Value = 8620
$ws.Range("K96").Value = 5185.7144
$ws.Range("L96").Value = 8620
$ws.Range("M96").Value = -3812.7144
$ws.Range("N96").Value = -11366

# WVR row 122: "Heavy Armoire" / "Dark Hempen Cloth"
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2320.6553
$ws.Range("I122").Value = 1409.4546
$ws.Range("J122").Value = 2877.5
$ws.Range("K122").Value = 4228.3638
$ws.Range("L122").Value = 8632.5
$ws.Range("M122").Value = -1778.3638
$ws.Range("N122").Value = -13532.5

# WVR row 126: "A Polished Purchase" / "Snow Linen"
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1568.84
$ws.Range("I126").Value = 1411.9412
$ws.Range("J126").Value = 1902.25
$ws.Range("K126").Value = 4235.8236
$ws.Range("L126").Value = 5706.75
$ws.Range("M126").Value = -1765.8236
$ws.Range("N126").Value = -10646.75
